$wb = $excel.ActiveWorkbook

# "Admin" sheet: D2 (Username) changes from "448924A" to "4482716A"
$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "4482716A"

# "Jira" sheet: B2 (Error description) changes from
# "0C6E5E6E11D4DE9764C9A5F1C9073D27" to "5FA3C2312892FD51F30690CB47131C4C"
$wsJira = $wb.Worksheets.Item("Jira")
$wsJira.Range("B2").Value = "5FA3C2312892FD51F30690CB47131C4C"
